$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the resistor values shown in B5/B6 (2.2K <-> 2.7K)
$ws.Range("B5").Value = "2.7K"
$ws.Range("B6").Value = "2.2K"

# Move the active selection from B28 to B9
$ws.Range("B9").Select()

# Reflect the updated application window position captured in the file
$excel.Left = 14460
$excel.Top = 80
